$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 493.16666
$ws.Range("J19").Value = 512.25
$ws.Range("L19").Value = 512.25
$ws.Range("N19").Value = -862.25
$ws.Range("H132").Value = 3543.2856
$ws.Range("I132").Value = 1241.8966
$ws.Range("K132").Value = 3725.6898
$ws.Range("M132").Value = -1195.6898
$ws.Range("H137").Value = 2420.1538
$ws.Range("I137").Value = 958.8570999999999
$ws.Range("J137").Value = 4125
$ws.Range("K137").Value = 2876.5713
$ws.Range("L137").Value = 12375
$ws.Range("M137").Value = -326.5712999999996
$ws.Range("N137").Value = -17475

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3671.7
$ws.Range("I74").Value = 2886.625
$ws.Range("K74").Value = 2886.625
$ws.Range("M74").Value = -2012.625
$ws.Range("H77").Value = 3671.7
$ws.Range("I77").Value = 2886.625
$ws.Range("K77").Value = 14433.125
$ws.Range("M77").Value = -10065.125
$ws.Range("H110").Value = 2527.75
$ws.Range("I110").Value = 2366
$ws.Range("J110").Value = 3013
$ws.Range("K110").Value = 2366
$ws.Range("L110").Value = 3013
$ws.Range("M110").Value = -321
$ws.Range("N110").Value = -7103

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 3302.75
$ws.Range("I75").Value = 3302.75
$ws.Range("K75").Value = 3302.75
$ws.Range("M75").Value = -2366.75
$ws.Range("H78").Value = 3302.75
$ws.Range("I78").Value = 3302.75
$ws.Range("K78").Value = 9908.25
$ws.Range("M78").Value = -5228.25
$ws.Range("H86").Value = 4348.4
$ws.Range("I86").Value = 1666.6666
$ws.Range("J86").Value = 5497.7144
$ws.Range("K86").Value = 1666.6666
$ws.Range("L86").Value = 5497.7144
$ws.Range("M86").Value = -543.6666
$ws.Range("N86").Value = -7743.7144
$ws.Range("H89").Value = 4348.4
$ws.Range("I89").Value = 1666.6666
$ws.Range("J89").Value = 5497.7144
$ws.Range("K89").Value = 8333.333000000001
$ws.Range("L89").Value = 27488.572
$ws.Range("M89").Value = -2717.333000000001
$ws.Range("N89").Value = -38720.572
$ws.Range("H107").Value = 9343.857
$ws.Range("I107").Value = 9402.833000000001
$ws.Range("J107").Value = 8990
$ws.Range("K107").Value = 9402.833000000001
$ws.Range("L107").Value = 8990
$ws.Range("M107").Value = -7482.833000000001
$ws.Range("N107").Value = -12830
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 5000101
$ws.Range("I3").Value = 5000101
$ws.Range("K3").Value = 5000101
$ws.Range("M3").Value = -4999988
$ws.Range("H5").Value = 1018.8461
$ws.Range("I5").Value = 468.14285
$ws.Range("J5").Value = 1661.3334
$ws.Range("K5").Value = 468.14285
$ws.Range("L5").Value = 1661.3334
$ws.Range("M5").Value = -356.14285
$ws.Range("N5").Value = -1885.3334
$ws.Range("H19").Value = 179
$ws.Range("I19").Value = 143.5
$ws.Range("J19").Value = 250
$ws.Range("K19").Value = 143.5
$ws.Range("L19").Value = 250
$ws.Range("M19").Value = 26.5
$ws.Range("N19").Value = -590
$ws.Range("H24").Value = 179
$ws.Range("I24").Value = 143.5
$ws.Range("J24").Value = 250
$ws.Range("K24").Value = 143.5
$ws.Range("L24").Value = 250
$ws.Range("M24").Value = 26.5
$ws.Range("N24").Value = -590
$ws.Range("H62").Value = 3005.2856
$ws.Range("I62").Value = 2957.8
$ws.Range("K62").Value = 2957.8
$ws.Range("M62").Value = -2333.8
$ws.Range("H65").Value = 3005.2856
$ws.Range("I65").Value = 2957.8
$ws.Range("K65").Value = 14789
$ws.Range("M65").Value = -11669

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 159.4
$ws.Range("I23").Value = 99
$ws.Range("J23").Value = 199.66667
$ws.Range("K23").Value = 297
$ws.Range("L23").Value = 599.00001
$ws.Range("M23").Value = -62
$ws.Range("N23").Value = -1069.00001
$ws.Range("H59").Value = 725
$ws.Range("I59").Value = 725
$ws.Range("K59").Value = 2175
$ws.Range("M59").Value = -1635
$ws.Range("H68").Value = 443
$ws.Range("I68").Value = 444
$ws.Range("J68").Value = 442
$ws.Range("K68").Value = 1332
$ws.Range("L68").Value = 1326
$ws.Range("M68").Value = -521
$ws.Range("N68").Value = -2948
$ws.Range("H71").Value = 443
$ws.Range("I71").Value = 444
$ws.Range("J71").Value = 442
$ws.Range("K71").Value = 3996
$ws.Range("L71").Value = 3978
$ws.Range("M71").Value = 60
$ws.Range("N71").Value = -12090
$ws.Range("H129").Value = 2123.5833
$ws.Range("I129").Value = 1075.25
$ws.Range("J129").Value = 2647.75
$ws.Range("K129").Value = 3225.75
$ws.Range("L129").Value = 7943.25
$ws.Range("M129").Value = 1774.25
$ws.Range("N129").Value = -17943.25
$ws.Range("H131").Value = 1838.3334
$ws.Range("I131").Value = 1258
$ws.Range("J131").Value = 2999
$ws.Range("K131").Value = 3774
$ws.Range("L131").Value = 8997
$ws.Range("M131").Value = 1266
$ws.Range("N131").Value = -19077
$ws.Range("H137").Value = 4362.25
$ws.Range("I137").Value = 899
$ws.Range("J137").Value = 5516.6665
$ws.Range("K137").Value = 2697
$ws.Range("L137").Value = 16549.9995
$ws.Range("M137").Value = 2403
$ws.Range("N137").Value = -26749.9995

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 35000
$ws.Range("J105").Value = 35000
$ws.Range("L105").Value = 35000
$ws.Range("N105").Value = -41988
$ws.Range("H113").Value = 7185.4287
$ws.Range("I113").Value = 3466
$ws.Range("K113").Value = 3466
$ws.Range("M113").Value = -1296
$ws.Range("H122").Value = 3851.7693
$ws.Range("I122").Value = 3798.889
$ws.Range("J122").Value = 3970.75
$ws.Range("K122").Value = 11396.667
$ws.Range("L122").Value = 11912.25
$ws.Range("M122").Value = -8946.667000000001
$ws.Range("N122").Value = -16812.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 825
$ws.Range("I31").Value = 825
$ws.Range("K31").Value = 825
$ws.Range("M31").Value = -577
$ws.Range("H100").Value = 7249.875
$ws.Range("I100").Value = 2666.3333
$ws.Range("K100").Value = 2666.3333
$ws.Range("M100").Value = -2125.3333

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 5000
$ws.Range("J3").Value = 5000
$ws.Range("L3").Value = 5000
$ws.Range("N3").Value = -5228
$ws.Range("H62").Value = 9958.25
$ws.Range("I62").Value = 8500.5
$ws.Range("J62").Value = 10687.125
$ws.Range("K62").Value = 8500.5
$ws.Range("L62").Value = 10687.125
$ws.Range("M62").Value = -7876.5
$ws.Range("N62").Value = -11935.125
$ws.Range("H65").Value = 9958.25
$ws.Range("I65").Value = 8500.5
$ws.Range("J65").Value = 10687.125
$ws.Range("K65").Value = 42502.5
$ws.Range("L65").Value = 53435.625
$ws.Range("M65").Value = -39382.5
$ws.Range("N65").Value = -59675.625
$ws.Range("H100").Value = 1498.5
$ws.Range("I100").Value = 1498.5
$ws.Range("K100").Value = 2997
$ws.Range("M100").Value = -2456
$ws.Range("H128").Value = 20000
$ws.Range("J128").Value = 20000
$ws.Range("L128").Value = 20000
$ws.Range("N128").Value = -29960
$ws.Range("H137").Value = 40079.5
$ws.Range("J137").Value = 40079.5
$ws.Range("L137").Value = 40079.5
$ws.Range("N137").Value = -50279.5
$ws.Range("H139").Value = 74900
$ws.Range("J139").Value = 74900
$ws.Range("L139").Value = 74900
$ws.Range("N139").Value = -85180
